$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Paragraph 2 reads "Version: 4.21" + <br/> + "Date: 10/6/2020" and the
# hidden "_GoBack" bookmark currently sits right after the version number
# (that was the last place Word's editor touched). This change bumps the
# version's patch digit (1 -> 3) and the date's day-of-month (6 -> 23);
# since the date edit is now the *last* text touched, "_GoBack" moves to
# the very end of the paragraph (right after "/2020", before the
# paragraph mark).
#
# NOTE: Range objects captured before an edit do not track later document
# mutations, so every position below is (re)located fresh, right before
# it's used.
# -------------------------------------------------------------------------

# --- 1. Version: "4.21" -> "4.23" (only the trailing "1" actually changes) ---
$verFind = $d.Content
$verFind.Find.Execute("4.21", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$verEnd = $verFind.End
$verDigit = $d.Range($verEnd - 1, $verEnd)   # the "1"
$verDigit.Text = "3"

# --- 2. Date: "10/6/2020" -> "10/23/2020" (only the "6" actually changes) ---
$dateFind = $d.Content
$dateFind.Find.Execute("10/6/2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dayStart = $dateFind.Start + 3   # skip past "10/"
$dayDigit = $d.Range($dayStart, $dayStart + 1)   # the "6"
$dayDigit.Text = "23"

# --- 3. Move the "_GoBack" bookmark to the end of the paragraph (after the
#        date, before the paragraph mark). Adding a bookmark with a
#        collapsed range sitting exactly one position before a paragraph
#        mark relocates incorrectly, so nudge the target off that boundary
#        with a scratch character, add the bookmark there, then remove the
#        scratch character again. ---
$dateFind2 = $d.Content
$dateFind2.Find.Execute("Date: 10/23/2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraEnd = $dateFind2.End   # position right before the paragraph mark

$scratch = $d.Range($paraEnd, $paraEnd)
$scratch.InsertAfter("X")

$bmSpot = $d.Range($paraEnd, $paraEnd)
$d.Bookmarks.Add("_GoBack", $bmSpot)

$d.Range($paraEnd, $paraEnd + 1).Delete()   # remove the "X" scratch char
